$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title: "Hand-off Document" + ":" + " " + "Team CSI (Northeastern
#    University)" were split across four runs; collapse them into a single
#    run. A Find/Replace over the already-correct text re-writes the match
#    as one run while preserving the shared run formatting (sz/szCs 32).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Hand-off Document: Team CSI (Northeastern University)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Hand-off Document: Team CSI (Northeastern University)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Our Imag" + "e can be written to an SD card, when the BBB is turned on
#    the new image will be installed. " were split across two runs; merge
#    them into one run. A temporary marker is appended right after the
#    sentence (still merged into the same run-producing replace) so that we
#    have a safe, non-paragraph-boundary position to anchor the relocated
#    "_GoBack" bookmark -- adding a bookmark exactly at a raw end-of-
#    paragraph offset is unreliable, so we avoid that edge case entirely.
# ---------------------------------------------------------------------------
$sentence = "Our Image can be written to an SD card, when the BBB is turned on the new image will be installed. "
$marker = "ZZ_BOOKMARK_MARKER_ZZ"
$rng = $d.Content
$rng.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1,
                   $false, ($sentence + $marker), 2) | Out-Null

# Re-create the "_GoBack" bookmark (collapsed) right before the marker text,
# i.e. right after the now-merged sentence -- exactly where it should sit at
# the end of this paragraph.
$markerStart = $rng.Start + $sentence.Length
$bm = $d.Range($markerStart, $markerStart)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null

# Remove the temporary marker text again.
$d.Content.Find.Execute($marker, $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Insert a new paragraph right after the "Our Image..." paragraph with a
#    note about the cape, inheriting that paragraph's formatting.
# ---------------------------------------------------------------------------
$para = $d.Paragraphs(13)
$para.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs(14)
$newPara.Range.Text = "Please note that you cannot have the cape on while attempting to install the image! "

# ---------------------------------------------------------------------------
# 4. "Please conta" + "ct us at " (previously split apart by the old
#    "_GoBack" bookmark sitting between them) are merged back into a single
#    run. Since the bookmark was already relocated above, this Find/Replace
#    both merges the runs and discards the now-stale bookmark markers that
#    used to live between them.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Please contact us at ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Please contact us at ", 2) | Out-Null
